# Noise Score * 0.1, general score adjust weight, clarity using minus sign
#
# - old_noise_score (C) and dl_noise_score (H) are rescaled by a factor of 0.1.
# - Restoration scores (F: old_restoration_score, K: dl_restoration_score,
#   P: self_restoration_score) are reweighted; P now also uses a minus sign
#   (can go negative) to reflect the clarity-adjustment convention.
# - The "self_*" group (L: self_clarity_score, M: self_noise_score,
#   N: self_contrast_score, O: self_color_score, P: self_restoration_score)
#   has its weighting adjusted across the board.
#
# Apply the updated values cell-by-cell for rows 2-17 (Picture 1..16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.222223211119789
$ws.Range("F2").Value = 2.559317730011622
$ws.Range("H2").Value = 4.625373463539849
$ws.Range("K2").Value = 4.596231881963304
$ws.Range("L2").Value = 17.28446471291866
$ws.Range("M2").Value = 4.991295985400049
$ws.Range("N2").Value = 19.57371478280446
$ws.Range("O2").Value = 24.12717311370862
$ws.Range("P2").Value = -2.083714412855725
$ws.Range("C3").Value = 3.603150636177189
$ws.Range("F3").Value = 2.992422884558274
$ws.Range("H3").Value = 5.298765881831494
$ws.Range("K3").Value = 4.517979248005615
$ws.Range("L3").Value = 6.225969414108619
$ws.Range("M3").Value = 3.969817732951864
$ws.Range("N3").Value = 15.56790661812184
$ws.Range("O3").Value = 20.41216723425281
$ws.Range("P3").Value = 1.473464961520643
$ws.Range("C4").Value = 6.47233922056953
$ws.Range("F4").Value = -0.2557722152980224
$ws.Range("H4").Value = 5.130961586720532
$ws.Range("K4").Value = 5.500266919955192
$ws.Range("L4").Value = 66.21822899573351
$ws.Range("M4").Value = 9.497148630804267
$ws.Range("N4").Value = 37.24370443203261
$ws.Range("O4").Value = 28.16507224873779
$ws.Range("P4").Value = -19.07118807625439
$ws.Range("C5").Value = 6.146790816742586
$ws.Range("F5").Value = 3.935859561147907
$ws.Range("H5").Value = 7.580557933448844
$ws.Range("K5").Value = 7.828240953729322
$ws.Range("L5").Value = 13.26006412674463
$ws.Range("M5").Value = 6.52599831047601
$ws.Range("N5").Value = 25.59214867785863
$ws.Range("O5").Value = 24.36611026990594
$ws.Range("P5").Value = 0.2932156187216668
$ws.Range("C6").Value = 6.518309386226512
$ws.Range("F6").Value = 4.06377356150258
$ws.Range("H6").Value = 4.270193039325234
$ws.Range("K6").Value = 3.477546967920692
$ws.Range("L6").Value = 10.03842026825633
$ws.Range("M6").Value = 6.777379886930184
$ws.Range("N6").Value = 26.57795115592092
$ws.Range("O6").Value = 19.77648636999793
$ws.Range("P6").Value = 1.244656794802388
$ws.Range("C7").Value = 6.556708149125992
$ws.Range("F7").Value = 4.98294434730147
$ws.Range("H7").Value = 6.75560858982304
$ws.Range("K7").Value = 6.111580346291397
$ws.Range("L7").Value = 12.11808844011142
$ws.Range("M7").Value = 6.963023438158046
$ws.Range("N7").Value = 27.30597826817089
$ws.Range("O7").Value = 33.21571584114716
$ws.Range("P7").Value = 1.846624830256911
$ws.Range("C8").Value = 2.66072378715509
$ws.Range("F8").Value = 3.922256241133697
$ws.Range("H8").Value = 2.664165717176652
$ws.Range("K8").Value = 4.345754245547921
$ws.Range("L8").Value = 5.413915973616289
$ws.Range("M8").Value = 2.816777870994172
$ws.Range("N8").Value = 11.04618637067342
$ws.Range("O8").Value = 32.21762696552921
$ws.Range("P8").Value = 2.420400219942837
$ws.Range("C9").Value = 4.415755032557821
$ws.Range("F9").Value = 2.5258365427665
$ws.Range("H9").Value = 3.663880465358452
$ws.Range("K9").Value = 3.597753601695072
$ws.Range("L9").Value = 25.94161733615222
$ws.Range("M9").Value = 6.245329638807759
$ws.Range("N9").Value = 24.49149171334361
$ws.Range("O9").Value = 34.120233169054
$ws.Range("P9").Value = -3.939924166529095
$ws.Range("C10").Value = 5.1167085298048
$ws.Range("F10").Value = 4.398130873978936
$ws.Range("H10").Value = 6.49995604114765
$ws.Range("K10").Value = 4.490367665280484
$ws.Range("L10").Value = 14.87813073094053
$ws.Range("M10").Value = 5.600666246139735
$ws.Range("N10").Value = 21.96339552362443
$ws.Range("O10").Value = 36.26854107813476
$ws.Range("P10").Value = 0.3880810463202282
$ws.Range("C11").Value = 6.536881959765681
$ws.Range("F11").Value = 4.733322261258856
$ws.Range("H11").Value = 6.596221032500619
$ws.Range("K11").Value = 4.968641779211397
$ws.Range("L11").Value = 10.52265320132967
$ws.Range("M11").Value = 6.925087483223781
$ws.Range("N11").Value = 27.15720405473038
$ws.Range("O11").Value = 27.49027811796294
$ws.Range("P11").Value = 1.893881097243367
$ws.Range("C12").Value = 5.828319888439398
$ws.Range("F12").Value = 3.723469196193645
$ws.Range("H12").Value = 6.606916958428691
$ws.Range("K12").Value = 5.623378826912037
$ws.Range("L12").Value = 17.08750890567113
$ws.Range("M12").Value = 6.703110813527243
$ws.Range("N12").Value = 26.28671813211271
$ws.Range("O12").Value = 32.68543386540362
$ws.Range("P12").Value = -0.3200497933637232
$ws.Range("C13").Value = 4.60627233099522
$ws.Range("F13").Value = 2.118439399750958
$ws.Range("H13").Value = 6.199419495745683
$ws.Range("K13").Value = 5.063324454896417
$ws.Range("L13").Value = 24.56446153558612
$ws.Range("M13").Value = 5.729196073446793
$ws.Range("N13").Value = 22.46744359205555
$ws.Range("O13").Value = 34.6814946209151
$ws.Range("P13").Value = -3.582905255765869
$ws.Range("C14").Value = 4.121325032780346
$ws.Range("F14").Value = 2.772763356404518
$ws.Range("H14").Value = 4.153789391957205
$ws.Range("K14").Value = 4.29073964997883
$ws.Range("L14").Value = 17.42986386516169
$ws.Range("M14").Value = 5.273489849649542
$ws.Range("N14").Value = 20.68035584584088
$ws.Range("O14").Value = 27.88213341330816
$ws.Range("P14").Value = -1.629707990460545
$ws.Range("C15").Value = 5.247114828446552
$ws.Range("F15").Value = 2.928097551347136
$ws.Range("H15").Value = 6.641397232977872
$ws.Range("K15").Value = 4.916638731943547
$ws.Range("L15").Value = 21.26021532896533
$ws.Range("M15").Value = 6.291667134464745
$ws.Range("N15").Value = 24.67319773979009
$ws.Range("O15").Value = 29.00189139299194
$ws.Range("P15").Value = -2.556757584668345
$ws.Range("C16").Value = 7.180938274560448
$ws.Range("F16").Value = 5.800639419285525
$ws.Range("H16").Value = 6.706349672525379
$ws.Range("K16").Value = 4.995953983757638
$ws.Range("L16").Value = 18.88898423005566
$ws.Range("M16").Value = 7.732042190561645
$ws.Range("N16").Value = 30.32174623946256
$ws.Range("O16").Value = 46.06058367156609
$ws.Range("P16").Value = 0.7952012658583634
$ws.Range("C17").Value = 6.444960907057918
$ws.Range("F17").Value = 3.664384470167237
$ws.Range("H17").Value = 4.859696759510381
$ws.Range("K17").Value = 4.010659367503363
$ws.Range("L17").Value = 14.472166283084
$ws.Range("M17").Value = 6.87898557916335
$ws.Range("N17").Value = 26.97641662188543
$ws.Range("O17").Value = 21.35423951342748
$ws.Range("P17").Value = -0.3218549112627729
